$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: "Team:" label and "${team}" expression (mirrors the Team Name:/Team City:/etc. block)
$ws.Range("A7").Value = "Team:"
$ws.Range("C7").Value = '${team}'

# Row 12: "List:" label and "${numberList}" expression (just above the "Avg:" row)
$ws.Range("A12").Value = "List:"
$ws.Range("B12").Value = '${numberList}'
